# v3.0 update FCI 27/1/2023
# Adds a new date column (C) with new values, and reorders the rows so the
# fund rows sit above the summary ("avg"/"total") rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell C1 needs the same style (bold, bordered, centered) as B1.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$ws.Range("C1").Value = "13-01-2023"

# Final row layout (label, B-value, C-value) after reordering.
$rows = @(
    @{ Label = "Delta Acciones";   B = 6358.49;  C = 6403.9 },
    @{ Label = "Delta Select";     B = 1978.7;   C = 1072.21 },
    @{ Label = "Fima Acciones";    B = 4085.57;  C = 3829.88 },
    @{ Label = "Fima PB Acciones"; B = 13668.89; C = 13210.12 },
    @{ Label = "avg";              B = 6522.91;  C = 6129.03 },
    @{ Label = "total";            B = 26091.65; C = 24516.11 }
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row.Label
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
}
